$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row for account 005683532 / SYLVERSON with balance
# 10051.66, placed right before the existing row for account 004224011
# (i.e. becomes the new row 6, pushing the rest of the table down by one).
$ws.Rows(6).Insert()
$ws.Cells.Item(6, 1).Value = "'005683532"
$ws.Cells.Item(6, 2).Value = "SYLVERSON"
$ws.Cells.Item(6, 3).Value = 10051.66

# Remove the old SYLVERSON row (account 005683532, balance 51.66), which
# after the insertion above has shifted down from row 101 to row 102.
$ws.Rows(102).Delete()
